$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "RegisterPage_Validdata" worksheet right after
#    "StackPageContent" (and before "LoginPageContent").
# ---------------------------------------------------------------------------
$stackSheet = $wb.Worksheets.Item("StackPageContent")
$ws = $wb.Worksheets.Add($null, $stackSheet)
$ws.Name = "RegisterPage_Validdata"

# ---------------------------------------------------------------------------
# 2. Populate the register-page test-data table (A1:F9).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "scenario_type"
$ws.Range("A1").NumberFormat = "@"
$ws.Range("B1").Value = "username"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("C1").Value = "password"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("D1").Value = "confirmpassword"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").Value = "submission_method"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("F1").Value = "expected_message"
$ws.Range("F1").NumberFormat = "@"
$ws.Range("A2").Value = "Null value in  username for register"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("C2").Value = "Qaninja@123"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("D2").Value = "Qaninja@123"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").Value = "submits the register form"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").Value = "Please fill out this field."
$ws.Range("F2").NumberFormat = "@"
$ws.Range("A3").Value = "Null value in password for register"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("B3").Value = "qatest"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("D3").Value = "Qaninja@123"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").Value = "submits the register form"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("F3").Value = "Please fill out this field."
$ws.Range("F3").NumberFormat = "@"
$ws.Range("A4").Value = "Null value in  confirm password for register"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("B4").Value = "qatest"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("C4").Value = "Qaninja@123"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("E4").Value = "submits the register form"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("F4").Value = "Please fill out this field."
$ws.Range("F4").NumberFormat = "@"
$ws.Range("A5").Value = "with specialcharacter password for register"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("B5").Value = "testdata3"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("C5").Value = "@@@`t@@@"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("D5").Value = "submits the register form"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").Value = "Password does not match requirement"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("A6").Value = "with password less then eight characters for register"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("B6").Value = "testdata3"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("C6").Value = "A!1"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("D6").Value = "A!1"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").Value = "submits the register form"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("F6").Value = "Password does not match requirement"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("A7").Value = "with mismatch password for register"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("B7").Value = "qatest123"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("C7").Value = "Qaninja@123"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("D7").Value = "Qaninja@1"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").Value = "submits the register form"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("F7").Value = "password_mismatch:The two password fields didn’t match."
$ws.Range("F7").NumberFormat = "@"
$ws.Range("A8").Value = "with specialcharacter username for register"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("B8").Value = "!@!"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").Value = "Qaninja@123"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("D8").Value = "Qaninja@123"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").Value = "submits the register form"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("F8").Value = "Username must be alphanumeric"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("A9").Value = "valid_register"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("B9").Value = "Validrun04"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("C9").Value = "Qaninja@123"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("D9").Value = "Qaninja@123"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").Value = "submits the register form"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("F9").Value = "New Account Created. You are logged in as username"
$ws.Range("F9").NumberFormat = "@"

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 3. Restore "StackPageContent"'s scroll position to A1 (clears the stale
#    topLeftCell="B1" left over from the previous session) without leaving
#    it as the active sheet.
# ---------------------------------------------------------------------------
$stackSheet.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# ---------------------------------------------------------------------------
# 4. Make the new Register sheet the active tab/sheet again, with its own
#    remembered selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B17").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Reposition the workbook window, mirroring the author's saved layout.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 6195
$win.Top = 4305
$win.Width = 21600
$win.Height = 11175

Write-Host "Register page sheet added and populated"
